# Adds the "ODI Bowling Extra" worksheet (mirrors the existing
# "ODI Batting Extra" sheet, but for bowling: MAIDEN_OVERS /
# PERCENT_WICKETS_OF_ALL) as the new last sheet in the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet, placed after the current last sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# ---------------------------------------------------------------------
# 2. Header row - reuse the exact header formatting already used by the
#    sibling "ODI Batting Extra" sheet (bold, bordered, centered).
# ---------------------------------------------------------------------
$styleSource = $wb.Worksheets.Item("ODI Batting Extra").Range("A1")
$styleSource.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# ---------------------------------------------------------------------
# 3. Data rows. Every value (even the numeric-looking ones) is stored
#    as plain text in the source data, so force a text number format
#    before writing so values like "4225"/"0"/"10.00%" aren't silently
#    re-interpreted as numbers/percentages.
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A2:C21")
$dataRange.NumberFormat = "@"

$rows = @(
    @("4225", $null, $null),
    @("4232", "0", "10.00%"),
    @("4233", $null, $null),
    @("4239", "0", $null),
    @("4252", "1", "60.00%"),
    @("4346", $null, $null),
    @("4402", $null, $null),
    @("4406", $null, $null),
    @("4410", "0", $null),
    @("4636", "0", "20.00%"),
    @("4639", "1", "40.00%"),
    @("4642", $null, $null),
    @("4648", "0", "10.00%"),
    @("4649", "1", "10.00%"),
    @("4669", "0", "30.00%"),
    @("4673", $null, $null),
    @("4676", "1", "20.00%"),
    @("4686", $null, $null),
    @("4688", "1", "20.00%"),
    @("4690", $null, $null)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $r = $r + 1
}

Write-Host "ODI Bowling Extra sheet added"
